# Decent progress on rush auto
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 23: add the same Encoder / tracking_enc device entry that already
# exists on row 22 (B22/C22) so it is repeated for row 23.
$ws.Range("B23").Value = "Encoder"
$ws.Range("C23").Value = "tracking_enc"

# Update the active view: scroll so A5 is the top-left cell and select D23.
$excel.ActiveWindow.TopLeftCell = $ws.Range("A5")
$ws.Range("D23").Select()
